$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" footer date field from
#    13/07/2022 -> 30/03/2023 everywhere it appears: the slide master and
#    every one of its custom (slide) layouts.
# ---------------------------------------------------------------------------
function Update-DateShapes {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "13/07/2022") {
                $tr.Text = "30/03/2023"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Figure 1.2 "Equivalent hours (kWh/kW):" label -> "Peak Sun Hours (kWh/kW):"
#    The single run is split into three runs ("Peak ", "Sun Hours ",
#    "(kWh/kW):") all sharing the same Arial / 14pt formatting.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Equivalent hours (kWh/kW):") {

            # Replace the text while keeping the existing run formatting
            # (Arial, 14pt) that was already applied to the whole range.
            $tr.Text = "Peak Sun Hours (kWh/kW):"

            # Force the text to split into three distinct runs matching the
            # three logical phrases, re-asserting the same font on each.
            $run1 = $tr.Characters(1, 5)    # "Peak "
            $run2 = $tr.Characters(6, 10)   # "Sun Hours "
            $run3 = $tr.Characters(16, 9)   # "(kWh/kW):"

            $run1.Font.Name = "Arial"
            $run1.Font.Size = 14

            $run2.Font.Name = "Arial"
            $run2.Font.Size = 14

            $run3.Font.Name = "Arial"
            $run3.Font.Size = 14
        }
    }
}
